$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.523.75'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '''3.082.19'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''543.12'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').Value = '''140.07'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''3.075.10'
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').Value = '''0.504'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '''0.458'
$ws.Range('D13').Value = '''35.08'
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = '''3.583.58'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '''63.566.04'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '''3.081.96'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '''6.66'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').Value = '''474.94'
$ws.Range('E20').Value = '  -3.72%  '
$ws.Range('D21').Value = '''13.47'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('E22').Value = '  -3.42%  '
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('D24').Value = '''78.86'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').Value = '''7.99'
$ws.Range('E28').Value = '  -6.19%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '''26.22'
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('E31').Value = '  -4.19%  '
$ws.Range('E32').Value = '  +2.03%  '
$ws.Range('D33').Value = '''58.12'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '''2.33'
$ws.Range('E34').Value = '  -7.82%  '
$ws.Range('D35').Value = '''5.46'
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('D36').Value = '''493.53'
$ws.Range('E36').Value = '  -5.62%  '
$ws.Range('D37').Value = '''6.01'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').Value = '''3.244.09'
$ws.Range('E38').Value = '  +2.84%  '
$ws.Range('D39').Value = '''0.0404'
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('D40').Value = '''0.0800'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').Value = '''2.63'
$ws.Range('E43').Value = '  -2.74%  '
$ws.Range('D44').Value = '''0.255'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '''25.67'
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('D47').Value = '''124.52'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').Value = '''2.05'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').Value = '''0.0₃0530'
$ws.Range('E49').Value = '  +3.50%  '
$ws.Range('D50').Value = '''0.110'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('E51').Value = '  +1.43%  '
